$wb = $excel.ActiveWorkbook

# --- Sheet "Informações Gerais" (1st sheet): insert 5 new case rows at row 39 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A39:H43").Insert(-4121)

$ws1.Range("A39").Value = 45102.08596229166
$ws1.Range("B39").Value = "0346.9/2023"
$ws1.Range("C39").Value = "Polícia Militar"
$ws1.Range("D39").Value = "Externo"
$ws1.Range("E39").Value = "Sgt. Jonathan"
$ws1.Range("F39").Value = 1031813
$ws1.Range("G39").Value = 11
$ws1.Range("H39").Value = "Apenas vítima(s)"

$ws1.Range("A40").Value = 45106.55609578703
$ws1.Range("B40").Value = "0518.9/2023"
$ws1.Range("C40").Value = "Polícia Militar"
$ws1.Range("D40").Value = "Externo"
$ws1.Range("E40").Value = "Sgt. Michel"
$ws1.Range("F40").Value = 9808655
$ws1.Range("G40").Value = 18
$ws1.Range("H40").Value = "Apenas vítima(s)"

$ws1.Range("A41").Value = 45107.06797804398
$ws1.Range("B41").Value = "0349.9/2023"
$ws1.Range("C41").Value = "Polícia Militar"
$ws1.Range("D41").Value = "Externo"
$ws1.Range("E41").Value = "St. Andreyer"
$ws1.Range("F41").Value = 9901221
$ws1.Range("G41").Value = 1
$ws1.Range("H41").Value = "Vítima(s) e veículo(s)"

$ws1.Range("A42").Value = 45112.384029444445
$ws1.Range("B42").Value = "0362.9/2023"
$ws1.Range("C42").Value = "Polícia Militar"
$ws1.Range("D42").Value = "Interno"
$ws1.Range("E42").Value = "Sgt. Andrelino"
$ws1.Range("F42").Value = 9901280
$ws1.Range("G42").Value = 11
$ws1.Range("H42").Value = "Apenas vítima(s)"

$ws1.Range("A43").Value = 45116.533558807874
$ws1.Range("B43").Value = "0361.9/2023"
$ws1.Range("C43").Value = "Polícia Militar"
$ws1.Range("D43").Value = "Externo"
$ws1.Range("E43").Value = "Sgt. R. Soares"
$ws1.Range("F43").Value = 136512
$ws1.Range("G43").Value = 17
$ws1.Range("H43").Value = "Apenas vítima(s)"

# --- Sheet "Vítimas" (2nd sheet): insert 5 matching rows at row 40 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A40:C44").Insert(-4121)

$ws2.Range("A40").Value = 45102.08635333333
$ws2.Range("B40").Value = "0346.9/2023"
$ws2.Range("C40").Value = 136534

$ws2.Range("A41").Value = 45106.557955208336
$ws2.Range("B41").Value = "0518.9/2023"
$ws2.Range("C41").Value = 137792

$ws2.Range("A42").Value = 45107.06833642361
$ws2.Range("B42").Value = "0349.9/2023"
$ws2.Range("C42").Value = 136528

$ws2.Range("A43").Value = 45112.38436751158
$ws2.Range("B43").Value = "0362.9/2023"
$ws2.Range("C43").Value = 136507

$ws2.Range("A44").Value = 45116.533809120374
$ws2.Range("B44").Value = "0361.9/2023"
$ws2.Range("C44").Value = 136512
